$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# --- Weekly crime stats table updates (rows 15-31) ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1

$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2

$ws.Range("I15").Value = 8

$ws.Range("L15").Value = 300

$ws.Range("M15").Value = 700

$ws.Range("N15").Value = 166.666666666667

$ws.Range("C16").Value = 12

$ws.Range("E16").Value = 100

$ws.Range("F16").Value = 20

$ws.Range("G16").Value = 20

$ws.Range("H16").Value = 0

$ws.Range("I16").Value = 66

$ws.Range("J16").Value = 90

$ws.Range("K16").Value = -26.666666666666

$ws.Range("L16").Value = 46.666666666666

$ws.Range("M16").Value = -13.157894736842

$ws.Range("N16").Value = -84.019370460048

$ws.Range("C17").Value = 3

$ws.Range("D17").Value = 7

$ws.Range("E17").Value = -57.142857142857

$ws.Range("F17").Value = 21

$ws.Range("G17").Value = 22

$ws.Range("H17").Value = -4.545454545454

$ws.Range("I17").Value = 123

$ws.Range("J17").Value = 114

$ws.Range("K17").Value = 7.894736842105

$ws.Range("L17").Value = 38.202247191011

$ws.Range("M17").Value = 161.702127659574

$ws.Range("N17").Value = -12.765957446808

$ws.Range("C18").NumberFormat = "General"
$ws.Range("C18").Value = "0"

$ws.Range("D18").Value = 8

$ws.Range("E18").Value = -100

$ws.Range("F18").Value = 19

$ws.Range("G18").Value = 14

$ws.Range("H18").Value = 35.714285714285

$ws.Range("J18").Value = 66

$ws.Range("K18").Value = 18.181818181818

$ws.Range("L18").Value = -4.878048780487

$ws.Range("M18").Value = 129.411764705882

$ws.Range("N18").Value = -70.676691729323

$ws.Range("C19").Value = 15

$ws.Range("D19").Value = 17

$ws.Range("E19").Value = -11.764705882352

$ws.Range("F19").Value = 50

$ws.Range("G19").Value = 52

$ws.Range("H19").Value = -3.846153846153

$ws.Range("I19").Value = 241

$ws.Range("J19").Value = 249

$ws.Range("K19").Value = -3.212851405622

$ws.Range("L19").Value = -2.032520325203

$ws.Range("M19").Value = 60.666666666666

$ws.Range("N19").Value = -30.144927536231

$ws.Range("C20").Value = 1

$ws.Range("G20").Value = 6

$ws.Range("H20").Value = 50

$ws.Range("I20").Value = 29

$ws.Range("K20").Value = 31.818181818181

$ws.Range("L20").Value = 16

$ws.Range("M20").Value = 61.111111111111

$ws.Range("N20").Value = -88.353413654618

$ws.Range("C21").Value = 32

$ws.Range("D21").Value = 38

$ws.Range("E21").Value = -15.78947368421

$ws.Range("F21").Value = 121

$ws.Range("G21").Value = 114

$ws.Range("H21").Value = 6.140350877192

$ws.Range("I21").Value = 545

$ws.Range("J21").Value = 542

$ws.Range("K21").Value = 0.553505535055

$ws.Range("L21").Value = 11.224489795918

$ws.Range("M21").Value = 67.177914110429

$ws.Range("N21").Value = -61.646727656579

$ws.Range("C22").NumberFormat = "General"
$ws.Range("C22").Value = "0"

$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Value = "0"

$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Value = "***.*"

$ws.Range("G22").Value = 2

$ws.Range("H22").Value = 50

$ws.Range("L22").Value = -4.761904761904

$ws.Range("M22").Value = -23.076923076923

$ws.Range("C23").NumberFormat = "General"
$ws.Range("C23").Value = "0"

$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 2

$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E23").Value = -100

$ws.Range("F23").Value = 4

$ws.Range("H23").Value = 0

$ws.Range("I23").Value = 12

$ws.Range("J23").Value = 19

$ws.Range("K23").Value = -36.842105263157

$ws.Range("L23").Value = -42.857142857142

$ws.Range("M23").Value = 100

$ws.Range("C24").Value = 31

$ws.Range("D24").Value = 34

$ws.Range("E24").Value = -8.823529411764

$ws.Range("F24").Value = 132

$ws.Range("G24").Value = 163

$ws.Range("H24").Value = -19.018404907975

$ws.Range("I24").Value = 582

$ws.Range("J24").Value = 767

$ws.Range("K24").Value = -24.119947848761

$ws.Range("L24").Value = -18.373071528751

$ws.Range("M24").Value = 21.757322175732

$ws.Range("D25").Value = 38

$ws.Range("E25").Value = -28.947368421052

$ws.Range("F25").Value = 98

$ws.Range("G25").Value = 144

$ws.Range("H25").Value = -31.944444444444

$ws.Range("I25").Value = 423

$ws.Range("J25").Value = 696

$ws.Range("K25").Value = -39.224137931034

$ws.Range("L25").Value = -35.419847328244

$ws.Range("C26").Value = 12

$ws.Range("D26").Value = 8

$ws.Range("E26").Value = 50

$ws.Range("F26").Value = 39

$ws.Range("G26").Value = 40

$ws.Range("H26").Value = -2.5

$ws.Range("I26").Value = 160

$ws.Range("J26").Value = 169

$ws.Range("K26").Value = -5.325443786982

$ws.Range("L26").Value = 20.300751879699

$ws.Range("M26").Value = 9.58904109589

$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1

$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 2

$ws.Range("I27").Value = 11

$ws.Range("K27").Value = 266.666666666667

$ws.Range("L27").Value = 450

$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Value = "0"

$ws.Range("E28").Value = -100

$ws.Range("F28").Value = 9

$ws.Range("H28").Value = 28.571428571428

$ws.Range("I28").Value = 33

$ws.Range("J28").Value = 19

$ws.Range("K28").Value = 73.684210526315

$ws.Range("L28").Value = 94.117647058823

$ws.Range("C31").NumberFormat = "General"
$ws.Range("C31").Value = "0"

$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Value = "0"

$ws.Range("E31").NumberFormat = "General"
$ws.Range("E31").Value = "***.*"

